# Applies the cryptocurrency price/volume update described by the commit
# "Updated cryptos list on Wed Dec 20 12:17:04 UTC 2023 with GitHub Actions".
# Every value is written as literal text (NumberFormat "@" while writing,
# then reset to the workbook's default "Normal" style) so that strings
# such as "42.10", "2.544.78" or "  -0.55%  " survive byte-for-byte instead
# of being auto-coerced into numbers/dates and losing padding or trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "42.796.37"
Set-TextValue "E2" "  -0.55%  "
Set-TextValue "D3" "2.206.78"
Set-TextValue "E3" "  -1.42%  "
Set-TextValue "E4" "  +0.05%  "
Set-TextValue "D5" "255.64"
Set-TextValue "E5" "  +0.83%  "
Set-TextValue "B6" "XRP"
Set-TextValue "C6" "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue "D6" "0.617"
Set-TextValue "E6" "  +1.20%  "
Set-TextValue "B7" "Solana"
Set-TextValue "C7" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D7" "76.82"
Set-TextValue "E7" "  +2.32%  "
Set-TextValue "E8" "  +0.08%  "
Set-TextValue "D9" "0.593"
Set-TextValue "E9" "  -1.19%  "
Set-TextValue "D10" "42.46"
Set-TextValue "E10" "  +1.77%  "
Set-TextValue "D11" "0.0908"
Set-TextValue "E11" "  -2.13%  "
Set-TextValue "D12" "6.99"
Set-TextValue "E12" "  +1.11%  "
Set-TextValue "E13" "  +1.04%  "
Set-TextValue "D14" "2.542.40"
Set-TextValue "E14" "  -1.28%  "
Set-TextValue "D15" "14.39"
Set-TextValue "E15" "  -1.49%  "
Set-TextValue "D16" "2.208.11"
Set-TextValue "E16" "  -1.68%  "
Set-TextValue "E17" "  -1.17%  "
Set-TextValue "D18" "42.771.48"
Set-TextValue "E18" "  -0.35%  "
Set-TextValue "E19" "  -0.86%  "
Set-TextValue "D20" "71.06"
Set-TextValue "E20" "  -0.27%  "
Set-TextValue "D21" "5.97"
Set-TextValue "E21" "  -0.41%  "
Set-TextValue "D22" "2.34"
Set-TextValue "E22" "  +7.32%  "
Set-TextValue "D23" "229.45"
Set-TextValue "E23" "  +0.42%  "
Set-TextValue "D24" "9.16"
Set-TextValue "E24" "  -5.46%  "
Set-TextValue "D26" "42.10"
Set-TextValue "E26" "  +6.50%  "
Set-TextValue "D27" "10.68"
Set-TextValue "E27" "  -0.68%  "
Set-TextValue "D28" "3.34"
Set-TextValue "E28" "  -3.41%  "
Set-TextValue "D29" "2.19"
Set-TextValue "E29" "  -2.20%  "
Set-TextValue "E30" "  +2.34%  "
Set-TextValue "D31" "172.20"
Set-TextValue "E31" "  +0.20%  "
Set-TextValue "D32" "20.30"
Set-TextValue "E32" "  +0.33%  "
Set-TextValue "D33" "0.0865"
Set-TextValue "E33" "  +8.06%  "
Set-TextValue "D34" "5.18"
Set-TextValue "E34" "  -1.91%  "
Set-TextValue "D35" "0.121"
Set-TextValue "E35" "  -0.29%  "
Set-TextValue "D36" "0.0355"
Set-TextValue "E36" "  +6.84%  "
Set-TextValue "D37" "0.106"
Set-TextValue "E37" "  -3.54%  "
Set-TextValue "D38" "4.33"
Set-TextValue "E38" "  -3.30%  "
Set-TextValue "D39" "13.05"
Set-TextValue "E39" "  +1.30%  "
Set-TextValue "D40" "2.89"
Set-TextValue "E40" "  +18.20%  "
Set-TextValue "D41" "2.10"
Set-TextValue "E41" "  -0.28%  "
Set-TextValue "B42" "MultiversX"
Set-TextValue "C42" "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue "D42" "61.54"
Set-TextValue "E42" "  +2.93%  "
Set-TextValue "B43" "Algorand"
Set-TextValue "C43" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D43" "0.201"
Set-TextValue "E43" "  -2.62%  "
Set-TextValue "D44" "5.26"
Set-TextValue "E44" "  -2.92%  "
Set-TextValue "D45" "102.68"
Set-TextValue "E45" "  -0.93%  "
Set-TextValue "D46" "8.46"
Set-TextValue "E46" "  -2.58%  "
Set-TextValue "B47" "Cronos"
Set-TextValue "C47" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D47" "0.0971"
Set-TextValue "E47" "  -1.73%  "
Set-TextValue "B48" "WOONetwork"
Set-TextValue "C48" "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
Set-TextValue "D48" "0.463"
Set-TextValue "E48" "  -4.08%  "
Set-TextValue "D49" "1.11"
Set-TextValue "E49" "  +0.07%  "
Set-TextValue "D50" "1.13"
Set-TextValue "E50" "  -1.51%  "
Set-TextValue "D51" "1.47"
Set-TextValue "E51" "  +22.03%  "
